# Adds the new "Algoritmo propuesto para revertir transformaciones
# (Ingenieria Inversa)" section (one bold heading paragraph followed by
# three body paragraphs) right after the last, empty paragraph that
# precedes the final <w:sectPr>.

$d = $word.ActiveDocument

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# Run/paragraph-mark properties shared by the whole section: Times New
# Roman, szCs 24, es-CO language -- bold variant adds <w:b/><w:bCs/>.
$rprBold = '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:bCs/><w:szCs w:val="24"/><w:lang w:val="es-CO"/></w:rPr>'
$rprNormal = '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:szCs w:val="24"/><w:lang w:val="es-CO"/></w:rPr>'

function Add-Paragraph([string]$innerXml) {
    # Appends a brand new paragraph (built from raw OOXML so the paragraph
    # mark's run-properties -- e.g. bCs -- come out exactly as intended)
    # right after the current last paragraph of the document body.
    $lastPara = $d.Paragraphs.Last
    [void]$lastPara.Range.InsertParagraphAfter()
    $newLast = $d.Paragraphs.Last
    $insertionRange = $newLast.Range
    [void]$insertionRange.Collapse(0)
    $fullXml = '<w:p ' + $wNs + '>' + $innerXml + '</w:p>'
    [void]$insertionRange.InsertXML($fullXml)
}

# 1) Bold section heading.
$heading = '<w:pPr>' + $rprBold + '</w:pPr>' + `
    '<w:r>' + $rprBold + '<w:t>Algoritmo propuesto para revertir transformaciones (Ingenier' + [char]0x00ED + 'a Inversa)</w:t></w:r>'
Add-Paragraph $heading

# 2) First body paragraph.
$body1Text = 'Para identificar qu' + [char]0x00E9 + ' transformaci' + [char]0x00F3 + 'n fue aplicada entre dos im' + [char]0x00E1 + 'genes consecutivas en el proceso de encriptado, se propone un algoritmo secuencial que prioriza operaciones reversibles. En primer lugar, se aplica la operaci' + [char]0x00F3 + 'n XOR entre la imagen encriptada y la imagen de referencia conocida como IM. Si esta verificaci' + [char]0x00F3 + 'n falla, se procede a probar con rotaciones de bits, ya que se ha descartado previamente el uso de desplazamientos por ser operaciones que implican p' + [char]0x00E9 + 'rdida irreversible de informaci' + [char]0x00F3 + 'n.'
$body1 = '<w:pPr>' + $rprNormal + '</w:pPr>' + `
    '<w:r>' + $rprNormal + '<w:t>' + $body1Text + '</w:t></w:r>'
Add-Paragraph $body1

# 3) Second body paragraph (trailing space -> xml:space="preserve").
$body2Text = 'Las rotaciones se prueban en ambas direcciones (izquierda y derecha), desde 1 hasta 7 bits. No se consideran rotaciones de 8 bits o superiores, ya que una rotaci' + [char]0x00F3 + 'n de 8 devuelve el mismo byte original, y valores mayores equivalen a rotaciones m' + [char]0x00E1 + 's peque' + [char]0x00F1 + 'as (por ejemplo, rotar 9 bits equivale a rotar 1). '
$body2 = '<w:pPr>' + $rprNormal + '</w:pPr>' + `
    '<w:r>' + $rprNormal + '<w:t xml:space="preserve">' + $body2Text + '</w:t></w:r>'
Add-Paragraph $body2

# 4) Third body paragraph, split around "txt" with spell-check markers,
#    same as the original author's text (".txt" file extension flagged by
#    the proofer).
$body3TextA = 'Cada posible transformaci' + [char]0x00F3 + 'n es validada mediante el proceso de enmascaramiento, que consiste en sumar los valores RGB resultantes con una m' + [char]0x00E1 + 'scara y comparar contra los valores registrados en los archivos .'
$body3TextC = ' desde la posici' + [char]0x00F3 + 'n indicada por la semilla.'
$body3 = '<w:pPr>' + $rprNormal + '</w:pPr>' + `
    '<w:r>' + $rprNormal + '<w:t>' + $body3TextA + '</w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r>' + $rprNormal + '<w:t>txt</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r>' + $rprNormal + '<w:t xml:space="preserve">' + $body3TextC + '</w:t></w:r>'
Add-Paragraph $body3

Write-Output ("Paragraphs after edit: " + $d.Paragraphs.Count)
